# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# to the latest scraped values, per the GitHub Actions refresh commit.
#
# Numeric-looking "Price" values (column D) are written with a leading
# apostrophe so Excel stores them as text (matching the source data's
# dotted/locale formatting, e.g. "1.012") instead of re-parsing them as
# numbers -- which would also silently drop meaningful trailing zeros
# (e.g. "0.3690" -> 0.369).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is: cell address, new value
$updates = @(
    @("D2", "27.440.14"),
    @("E2", "  +1.97%  "),
    @("D3", "1.835.34"),
    @("E3", "  +1.09%  "),
    @("D4", "'1.012"),
    @("E4", "  +1.04%  "),
    @("D5", "'314.27"),
    @("E5", "  +1.61%  "),
    @("E6", "  +0.88%  "),
    @("D7", "'0.4739"),
    @("E7", "  +2.04%  "),
    @("D8", "'0.3690"),
    @("E8", "  +0.77%  "),
    @("D9", "'0.07456"),
    @("E9", "  +1.39%  "),
    @("D10", "'0.8848"),
    @("E10", "  +1.72%  "),
    @("D11", "'20.47"),
    @("E11", "  +0.83%  "),
    @("D12", "1.879.35"),
    @("E12", "  +2.54%  "),
    @("D13", "'0.07348"),
    @("E13", "  +3.70%  "),
    @("D14", "'5.447"),
    @("E14", "  +1.35%  "),
    @("D15", "'93.15"),
    @("E15", "  +1.85%  "),
    @("D16", "'6.579"),
    @("E16", "  +1.13%  "),
    @("D17", "'1.012"),
    @("E17", "  +0.89%  "),
    @("E18", "  +1.21%  "),
    @("E19", "  +0.94%  "),
    @("B20", "WrappedBTC"),
    @("C20", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"),
    @("D20", "27.606.73"),
    @("E20", "  +2.46%  "),
    @("B21", "Avalanche"),
    @("C21", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"),
    @("D21", "'14.80"),
    @("E21", "  +1.04%  "),
    @("D22", "'5.314"),
    @("E22", "  +0.28%  "),
    @("D23", "'10.68"),
    @("E23", "  +0.44%  "),
    @("D24", "2.102.13"),
    @("E24", "  +2.58%  "),
    @("D25", "'1.910"),
    @("E25", "  +0.78%  "),
    @("D26", "'152.11"),
    @("E26", "  +0.83%  "),
    @("D27", "'18.64"),
    @("E27", "  +1.79%  "),
    @("E28", "  +0.56%  "),
    @("D29", "'5.244"),
    @("E29", "  -0.22%  "),
    @("D30", "'117.67"),
    @("E30", "  +2.20%  "),
    @("D31", "'0.08997"),
    @("E31", "  +1.08%  "),
    @("D32", "'0.7561"),
    @("E32", "  -0.16%  "),
    @("E33", "  +1.88%  "),
    @("D34", "'4.549"),
    @("E34", "  +1.48%  "),
    @("D35", "'2.948"),
    @("E35", "  +1.29%  "),
    @("E36", "  +1.04%  "),
    @("D37", "'1.103"),
    @("E37", "  +1.68%  "),
    @("D38", "'0.05336"),
    @("E38", "  +1.13%  "),
    @("D39", "'0.01955"),
    @("E39", "  +0.50%  "),
    @("D40", "'2.984"),
    @("E40", "  +0.05%  "),
    @("D41", "'7.314"),
    @("E41", "  +1.17%  "),
    @("D42", "'2.407"),
    @("E42", "  +5.34%  "),
    @("D43", "'0.5326"),
    @("E43", "  +0.59%  "),
    @("D44", "'0.1659"),
    @("E44", "  +0.30%  "),
    @("D45", "'8.512"),
    @("E45", "  +1.06%  "),
    @("D46", "'0.4905"),
    @("E46", "  +0.76%  "),
    @("D47", "'10.52"),
    @("E47", "  +0.66%  "),
    @("E48", "  +1.00%  "),
    @("D49", "'104.91"),
    @("E49", "  +1.58%  "),
    @("D50", "'1.676"),
    @("E50", "  +0.97%  "),
    @("D51", "'0.06304"),
    @("E51", "  +0.19%  "),
)

foreach ($update in $updates) {
    $cell = $update[0]
    $value = $update[1]
    $ws.Range($cell).Value = $value
}
